# "Update countries & provincias Spain"
# Refresh the COVID-19 "Pais" sheet with the newer snapshot (24 Aug 2020,
# 18:36 instead of 17:19): updated timestamp, updated per-country metrics,
# and a handful of countries that changed rank (so the same table rows now
# correspond to different countries, with freshly updated figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Title / timestamp banner (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 24 de Agosto de 2020 a las 18:36"

# --- Per-row data refresh (B=Casos totales, C=Nuevos casos, D=Casos activos,
#     E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 5887116
$ws.Cells.Item(4, 3).Value = 12970
$ws.Cells.Item(4, 4).Value = 3169897
$ws.Cells.Item(4, 5).Value = 2536491
$ws.Cells.Item(4, 7).Value = 124
$ws.Cells.Item(4, 8).Value = 180728

# Row 6: India
$ws.Cells.Item(6, 2).Value = 3149759
$ws.Cells.Item(6, 3).Value = 44574
$ws.Cells.Item(6, 5).Value = 735619
$ws.Cells.Item(6, 7).Value = 625
$ws.Cells.Item(6, 8).Value = 58317

# Row 16: Reino Unido
$ws.Cells.Item(16, 7).Value = 4
$ws.Cells.Item(16, 8).Value = 41433

# Row 20: Italia
$ws.Cells.Item(20, 2).Value = 260298
$ws.Cells.Item(20, 3).Value = 953
$ws.Cells.Item(20, 4).Value = 205662
$ws.Cells.Item(20, 5).Value = 19195
$ws.Cells.Item(20, 7).Value = 4
$ws.Cells.Item(20, 8).Value = 35441

# Row 23: Alemania
$ws.Cells.Item(23, 2).Value = 235614
$ws.Cells.Item(23, 3).Value = 1125
$ws.Cells.Item(23, 5).Value = 16680

# Row 49: Singapur
$ws.Cells.Item(49, 4).Value = 54587
$ws.Cells.Item(49, 5).Value = 1790

# Row 74: Chequia
$ws.Cells.Item(74, 2).Value = 22056
$ws.Cells.Item(74, 3).Value = 133
$ws.Cells.Item(74, 4).Value = 16342
$ws.Cells.Item(74, 5).Value = 5299
$ws.Cells.Item(74, 7).Value = 3
$ws.Cells.Item(74, 8).Value = 415

# Rows 85-87: Libano / Senegal / Sudan move up one rank each (Libano
# overtakes Senegal and Sudan), with refreshed figures for the new snapshot.
$ws.Cells.Item(85, 1).Value = "Libano"
$ws.Cells.Item(85, 2).Value = 13155
$ws.Cells.Item(85, 3).Value = 457
$ws.Cells.Item(85, 4).Value = 3704
$ws.Cells.Item(85, 5).Value = 9325
$ws.Cells.Item(85, 8).Value = 126

$ws.Cells.Item(86, 1).Value = "Senegal"
$ws.Cells.Item(86, 2).Value = 13013
$ws.Cells.Item(86, 3).Value = 64
$ws.Cells.Item(86, 4).Value = 8595
$ws.Cells.Item(86, 5).Value = 4146
$ws.Cells.Item(86, 7).Value = 3
$ws.Cells.Item(86, 8).Value = 272

$ws.Cells.Item(87, 1).Value = "Sudan"
$ws.Cells.Item(87, 2).Value = 12836
$ws.Cells.Item(87, 4).Value = 6497
$ws.Cells.Item(87, 5).Value = 5524
$ws.Cells.Item(87, 8).Value = 815

# Row 90: Noruega
$ws.Cells.Item(90, 2).Value = 10374
$ws.Cells.Item(90, 3).Value = 51
$ws.Cells.Item(90, 5).Value = 960

# Rows 94-95: Grecia overtakes Guayana Francesa.
$ws.Cells.Item(94, 1).Value = "Grecia"
$ws.Cells.Item(94, 2).Value = 8819
$ws.Cells.Item(94, 3).Value = 155
$ws.Cells.Item(94, 4).Value = 3804
$ws.Cells.Item(94, 5).Value = 4773
$ws.Cells.Item(94, 8).Value = 242

$ws.Cells.Item(95, 1).Value = "Guayana Francesa"
$ws.Cells.Item(95, 2).Value = 8797
$ws.Cells.Item(95, 4).Value = 8307
$ws.Cells.Item(95, 5).Value = 435
$ws.Cells.Item(95, 8).Value = 55

# Row 96: Albania
$ws.Cells.Item(96, 2).Value = 8605
$ws.Cells.Item(96, 3).Value = 178
$ws.Cells.Item(96, 4).Value = 4413
$ws.Cells.Item(96, 5).Value = 3938
$ws.Cells.Item(96, 7).Value = 4
$ws.Cells.Item(96, 8).Value = 254

# Row 100: Haiti
$ws.Cells.Item(100, 2).Value = 8110
$ws.Cells.Item(100, 3).Value = 28
$ws.Cells.Item(100, 5).Value = 2290

# Row 104: Maldivas
$ws.Cells.Item(104, 5).Value = 2530
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = 27

# Row 128: Mali
$ws.Cells.Item(128, 2).Value = 2708
$ws.Cells.Item(128, 3).Value = 3
$ws.Cells.Item(128, 4).Value = 2025
$ws.Cells.Item(128, 5).Value = 558

# Row 145: Jordania
$ws.Cells.Item(145, 2).Value = 1639
$ws.Cells.Item(145, 3).Value = 30
$ws.Cells.Item(145, 4).Value = 1335
$ws.Cells.Item(145, 5).Value = 290
$ws.Cells.Item(145, 7).Value = 2
$ws.Cells.Item(145, 8).Value = 14

# Rows 214-215: Islas Malvinas and Montserrat swap order.
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0

$ws.Cells.Item(215, 1).Value = "Montserrat"
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 8).Value = 1
